# Applies the "Issues fixes and reports" commit:
#   - AMSIN sheet: append rows 74 and 75 (two new sprint runs)
#   - BETA sheet: append row 31 (one new sprint run)
#   - AMS sheet:  append row 41 (one new sprint run) and normalize the
#                 cell formatting on row 40 to match the rest of the table
#
# xlPasteFormats
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

function Add-SprintRow($ws, $row, $formatSourceRow, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    # Run Date / Sprint Name are literal text (not dates) - force text format
    # first so Excel doesn't reinterpret "2023-02-17" style strings as dates.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $sprintName

    # Total/Pass/Fail/Time Taken are plain numbers on a brand-new row, so
    # they naturally pick up the table's normal cell formatting.
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken

    # Run Time keeps the datetime number format used throughout column B.
    $bSrc = $ws.Range("B" + $formatSourceRow)
    $bSrc.Copy()
    $ws.Range("B" + $row).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 2).Value = $runTime
}

# Same as Add-SprintRow, but first clones the whole-row formatting from
# $formatSourceRow (used where the target row should look un-styled, i.e.
# match a still-unformatted source row's plain appearance).
function Add-SprintRowLikeSource($ws, $row, $formatSourceRow, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    $srcRange = $ws.Range("A" + $formatSourceRow + ":G" + $formatSourceRow)
    $srcRange.Copy()
    $dstRange = $ws.Range("A" + $row + ":G" + $row)
    $dstRange.PasteSpecial($xlPasteFormats)

    Add-SprintRow $ws $row $formatSourceRow $runDate $runTime $sprintName $total $pass $fail $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN sheet: dimension A1:G73 -> A1:G75
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-SprintRow $wsAmsin 74 73 "2023-02-17" 44974.44401239583 "173cyclefst" 89 87 2 2.08
Add-SprintRow $wsAmsin 75 73 "2023-02-20" 44977.41447871528 "173fnlrun" 89 87 2 2.23

# ---------------------------------------------------------------------
# BETA sheet: dimension A1:G30 -> A1:G31
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-SprintRow $wsBeta 31 30 "2023-02-20" 44977.59371335648 "173beta" 89 89 0 2.26

# ---------------------------------------------------------------------
# AMS sheet: dimension A1:G40 -> A1:G41
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 40 itself gets re-written first so its cells pick up the same
# formatting as the rest of the table (it previously had no explicit
# style). Deleting+re-entering it must happen before row 41 is added,
# otherwise deleting row 40 would wipe out the freshly added row 41.
$wsAms.Rows(40).Delete()
$wsAms.Cells.Item(40, 1).NumberFormat = "@"
$wsAms.Cells.Item(40, 1).Value = "2023-01-20"
$wsAms.Cells.Item(40, 3).NumberFormat = "@"
$wsAms.Cells.Item(40, 3).Value = "172live"
$wsAms.Cells.Item(40, 4).Value = 89
$wsAms.Cells.Item(40, 5).Value = 89
$wsAms.Cells.Item(40, 6).Value = 0
$wsAms.Cells.Item(40, 7).Value = 2.31

$wsAms.Range("B39").Copy()
$wsAms.Range("B40").PasteSpecial($xlPasteFormats)
$wsAms.Cells.Item(40, 2).Value = 44946.90048298611

# Row 41 is added after, copying row 40's freshly-normalized look; the
# numeric columns resolve back to the plain (un-styled) appearance that
# row 41 has in the target file since that format is visually a no-op.
Add-SprintRowLikeSource $wsAms 41 40 "2023-02-20" 44977.84407370233 "live173" 89 89 0 2.4
